# Scraper data + show data MVC
#
# A newer Marsa Maroc tender ("TRAVAUX DE VIDANGE DES BACS DE STOCKAGE ...
# MOHAMMEDIA", due Jeu 10 Juil 2025) was picked up by the scraper. It is
# rendered at the top of the tenders table, every older tender shifts down
# one row, and the tender that falls off the bottom of the fixed-size table
# is replaced by the next one in line ("Travaux mécaniques, électriques et
# divers des engins flottants ... Port d'Agadir"), which keeps the date of
# the row it now occupies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tenders = @(
    @{ Objet = "Objet : TRAVAUX DE VIDANGE DES BACS DE STOCKAGE DE LA STATION DEBALLASTAGE ET RINÇAGE DE MARSA MAROC AU PORT DE MOHAMMEDIA"; DateLimite = "Jeu 10 Juil 2025" },
    @{ Objet = "Objet : Travaux de construction de murettes amovibles en béton armé pour le compte de Marsa Maroc et sa filiale SMA au Port d’Agadir"; DateLimite = "Mar 08 Juil 2025" },
    @{ Objet = "Objet : TRAVAUX D’AMENAGEMENT DU TERMINAL A CONTENEURS EST DU PORT DE NADOR WEST MED LOT : ELECTRIFICATION HT, BT ET ECT"; DateLimite = "Jeu 03 Juil 2025" },
    @{ Objet = "Objet : Fourniture et installation des coffrets électriques pour le compte de MarsaMaroc au port de Tanger Med 1"; DateLimite = "Mer 02 Juil 2025" },
    @{ Objet = "Objet : FOURNITURE DE CHARIOTS ELEVATEURS DE MOYENNE CAPACITE A MARSA MAROC"; DateLimite = "Mer 02 Juil 2025" },
    @{ Objet = "Objet : Démolition d’un ancien bâtiment relevant de Marsa Maroc à la Direction de l’Exploitation au port de Casablanca Trafic Polyvalent"; DateLimite = "Mar 01 Juil 2025" },
    @{ Objet = "Objet : La fourniture des remorques basses et des remorques à cuvettes au Terminal à Conteneurs Est de Marsa Maroc Trafics Conteneur et Roulier et à la société Terminal à Conteneurs 3 du Port de Casablanca."; DateLimite = "Mar 01 Juil 2025" },
    @{ Objet = "Objet : SOUS TRAITANCE DES PRESATIONS DE GERBAGE ET CHARGEMENT DES MARCHANDISES PAR LES EQUIPEMENTS DE MARSA MAROC AU PORT DE JORF LASFAR"; DateLimite = "Lun 30 Juin 2025" },
    @{ Objet = "Objet : construction de deux hangars pour le stockage des vracs solides au Terminal Polyvalent de Marsa Maroc au Port de Casablanca."; DateLimite = "Jeu 26 Juin 2025" },
    @{ Objet = "Objet : Travaux mécaniques, électriques et divers des engins flottants de la société Marsa Maroc au Port d’Agadir"; DateLimite = "Jeu 26 Juin 2025" }
)

$row = 2
foreach ($tender in $tenders) {
    $ws.Cells.Item($row, 1).Value = $tender.Objet
    $ws.Cells.Item($row, 2).Value = $tender.DateLimite
    $row = $row + 1
}
